$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.386.68'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '3.374.52'
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.56%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.374.64'
$ws.Range("E8").Value = '  -2.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.471'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.64%  '
$ws.Range("E11").Value = '  -3.10%  '
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("D13").Value = '3.951.08'
$ws.Range("E13").Value = '  -2.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D16").Value = '3.381.88'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000170'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.55%  '
$ws.Range("D18").Value = '60.539.32'
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("E19").Value = '  -1.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.19%  '
$ws.Range("E21").Value = '  -5.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '386.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("E23").Value = '  -2.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -8.20%  '
$ws.Range("D27").Value = '3.517.98'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.178'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.05%  '
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("E33").Value = '  -9.12%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.54%  '
$ws.Range("D36").Value = '3.405.30'
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0768'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -1.97%  '
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.38%  '
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D48").Value = '2.515.19'
$ws.Range("E49").Value = '  -4.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.07'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("E51").Value = '  -3.35%  '
